# Agents.xlsx restructuring:
#  - drop the "Sim2" sheet (workbook keeps only "Sim1")
#  - replace the 3 sample agent rows with 5 new ones (new names/strategies)
#  - widen column B to fit the new names
#  - leave the selection on Sim1!C7 (tab now selected there)

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# --- remove the Sim2 sheet -------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sim2")
[void]$ws2.Delete()

$ws = $wb.Worksheets.Item("Sim1")

# --- rewrite the agent table (rows 2-6) ------------------------------------
$ws.Cells.Item(2,1).Value = "producer"
$ws.Cells.Item(2,2).Value = "GasPlant1"
$ws.Cells.Item(2,3).Value = "NaturalGasBiddingStrategy"
$ws.Cells.Item(2,4).Value = "{}"

$ws.Cells.Item(3,1).Value = "producer"
$ws.Cells.Item(3,2).Value = "CoalPlant1"
$ws.Cells.Item(3,3).Value = "CoalBiddingStrategy"
$ws.Cells.Item(3,4).Value = "{}"

$ws.Cells.Item(4,1).Value = "producer"
$ws.Cells.Item(4,2).Value = "HydroPlant1"
$ws.Cells.Item(4,3).Value = "DammedHydroBiddingStrategy"
$ws.Cells.Item(4,4).Value = "{}"

$ws.Cells.Item(5,1).Value = "producer"
$ws.Cells.Item(5,2).Value = "ZeroBidders1"
$ws.Cells.Item(5,3).Value = "ZeroBiddingStrategy"
$ws.Cells.Item(5,4).Value = "{}"

$ws.Cells.Item(6,1).Value = "consumer"
$ws.Cells.Item(6,2).Value = "Consumer"
$ws.Cells.Item(6,3).Value = "ConsumerBiddingStrategy"
$ws.Cells.Item(6,4).Value = "{}"

# --- column B needs to fit the longer agent names --------------------------
$ws.Columns.Item(2).ColumnWidth = 11.1

# --- restore selection/active sheet state -----------------------------------
[void]$ws.Activate()
[void]$ws.Range("C7").Select()
